# Refresh the cryptos price/volume(1h) snapshot (GitHub Actions scrape update),
# plus a ranking reshuffle among rows 37-40 (Dai / PEPE / InjectiveProtocol / Maker).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price/volume columns are stored as plain TEXT (e.g. "0.999", "65.567.42") even when
# the text happens to look like a number. Excel auto-converts a numeric-looking string
# to a real number on assignment unless the cell is already formatted as Text, so force
# that only where needed, then drop back to the Normal style to avoid leaving stray
# number-formatting behind on cells that do not otherwise have any.
function Set-TextCell($range, [string]$value) {
    if ($value -match '^[+-]?[0-9]*\.?[0-9]+([eE][+-]?[0-9]+)?$') {
        $range.NumberFormat = "@"
        $range.Value = $value
        $range.Style = "Normal"
    } else {
        $range.Value = $value
    }
}

Set-TextCell $ws.Range("D2") '65.567.42'
Set-TextCell $ws.Range("E2") '  -3.20%  '
Set-TextCell $ws.Range("D3") '3.504.87'
Set-TextCell $ws.Range("E3") '  +0.00%  '
Set-TextCell $ws.Range("D4") '0.999'
Set-TextCell $ws.Range("E4") '  +0.09%  '
Set-TextCell $ws.Range("D5") '556.72'
Set-TextCell $ws.Range("E5") '  +0.13%  '
Set-TextCell $ws.Range("D6") '179.61'
Set-TextCell $ws.Range("E6") '  -6.24%  '
Set-TextCell $ws.Range("D7") '0.641'
Set-TextCell $ws.Range("E7") '  +4.40%  '
Set-TextCell $ws.Range("E8") '  +0.11%  '
Set-TextCell $ws.Range("D9") '0.632'
Set-TextCell $ws.Range("E9") '  -1.37%  '
Set-TextCell $ws.Range("D10") '0.155'
Set-TextCell $ws.Range("E10") '  +2.93%  '
Set-TextCell $ws.Range("D11") '53.95'
Set-TextCell $ws.Range("E11") '  -6.35%  '
Set-TextCell $ws.Range("D12") '0.0000274'
Set-TextCell $ws.Range("E12") '  -1.11%  '
Set-TextCell $ws.Range("D13") '9.28'
Set-TextCell $ws.Range("E13") '  -2.45%  '
Set-TextCell $ws.Range("D14") '4.065.86'
Set-TextCell $ws.Range("E14") '  +0.01%  '
Set-TextCell $ws.Range("D15") '3.503.24'
Set-TextCell $ws.Range("E15") '  -0.10%  '
Set-TextCell $ws.Range("D16") '18.47'
Set-TextCell $ws.Range("E16") '  +0.37%  '
Set-TextCell $ws.Range("E17") '  +0.31%  '
Set-TextCell $ws.Range("D18") '12.15'
Set-TextCell $ws.Range("E18") '  +2.57%  '
Set-TextCell $ws.Range("D19") '65.555.57'
Set-TextCell $ws.Range("E19") '  -3.52%  '
Set-TextCell $ws.Range("D20") '0.998'
Set-TextCell $ws.Range("E20") '  -1.33%  '
Set-TextCell $ws.Range("E21") '  +2.50%  '
Set-TextCell $ws.Range("E22") '  +2.71%  '
Set-TextCell $ws.Range("D23") '86.14'
Set-TextCell $ws.Range("E23") '  +1.39%  '
Set-TextCell $ws.Range("D24") '4.27'
Set-TextCell $ws.Range("E24") '  +0.98%  '
Set-TextCell $ws.Range("D25") '12.85'
Set-TextCell $ws.Range("E25") '  +8.26%  '
Set-TextCell $ws.Range("D26") '10.83'
Set-TextCell $ws.Range("E26") '  -9.23%  '
Set-TextCell $ws.Range("E27") '  -2.08%  '
Set-TextCell $ws.Range("D28") '6.05'
Set-TextCell $ws.Range("E28") '  -3.57%  '
Set-TextCell $ws.Range("D29") '9.07'
Set-TextCell $ws.Range("E29") '  +4.78%  '
Set-TextCell $ws.Range("D30") '30.41'
Set-TextCell $ws.Range("E30") '  -0.35%  '
Set-TextCell $ws.Range("D31") '6.50'
Set-TextCell $ws.Range("E31") '  -5.70%  '
Set-TextCell $ws.Range("D32") '607.89'
Set-TextCell $ws.Range("E32") '  -12.23%  '
Set-TextCell $ws.Range("E33") '  -0.36%  '
Set-TextCell $ws.Range("E34") '  -0.89%  '
Set-TextCell $ws.Range("D35") '59.75'
Set-TextCell $ws.Range("E35") '  -1.71%  '
Set-TextCell $ws.Range("D36") '0.146'
Set-TextCell $ws.Range("E36") '  +9.13%  '
Set-TextCell $ws.Range("B37") 'Dai'
Set-TextCell $ws.Range("C37") 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextCell $ws.Range("D37") '1.00'
Set-TextCell $ws.Range("E37") '  +0.46%  '
Set-TextCell $ws.Range("B38") 'PEPE'
Set-TextCell $ws.Range("C38") 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
Set-TextCell $ws.Range("D38") '0.0₃0795'
Set-TextCell $ws.Range("E38") '  -4.05%  '
Set-TextCell $ws.Range("B39") 'InjectiveProtocol'
Set-TextCell $ws.Range("C39") 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextCell $ws.Range("D39") '37.39'
Set-TextCell $ws.Range("E39") '  -4.31%  '
Set-TextCell $ws.Range("B40") 'Maker'
Set-TextCell $ws.Range("C40") 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextCell $ws.Range("D40") '3.385.41'
Set-TextCell $ws.Range("E40") '  +10.79%  '
Set-TextCell $ws.Range("D41") '0.381'
Set-TextCell $ws.Range("E41") '  -5.78%  '
Set-TextCell $ws.Range("D42") '0.999'
Set-TextCell $ws.Range("E42") '  -0.29%  '
Set-TextCell $ws.Range("E43") '  -3.92%  '
Set-TextCell $ws.Range("D44") '2.84'
Set-TextCell $ws.Range("E44") '  -4.80%  '
Set-TextCell $ws.Range("E45") '  -9.34%  '
Set-TextCell $ws.Range("E46") '  -1.74%  '
Set-TextCell $ws.Range("D47") '3.24'
Set-TextCell $ws.Range("E47") '  +0.78%  '
Set-TextCell $ws.Range("D48") '2.70'
Set-TextCell $ws.Range("E48") '  -2.57%  '
Set-TextCell $ws.Range("E49") '  +1.79%  '
Set-TextCell $ws.Range("D50") '8.46'
Set-TextCell $ws.Range("E50") '  -4.51%  '
Set-TextCell $ws.Range("D51") '137.97'
Set-TextCell $ws.Range("E51") '  -1.89%  '
